$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-21 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-22 Thursday", 2) | Out-Null
$d.Content.Find.Execute("34×88=2992", $true, $false, $false, $false, $false, $true, 1, $false, "97×90=8730", 2) | Out-Null
$d.Content.Find.Execute("37×42=1554", $true, $false, $false, $false, $false, $true, 1, $false, "42×35=1470", 2) | Out-Null
$d.Content.Find.Execute("38×19=722", $true, $false, $false, $false, $false, $true, 1, $false, "94×24=2256", 2) | Out-Null
$d.Content.Find.Execute("87×43=3741", $true, $false, $false, $false, $false, $true, 1, $false, "76×30=2280", 2) | Out-Null
$d.Content.Find.Execute("32×73=2336", $true, $false, $false, $false, $false, $true, 1, $false, "62×94=5828", 2) | Out-Null
$d.Content.Find.Execute("42×40=1680", $true, $false, $false, $false, $false, $true, 1, $false, "72×10=720", 2) | Out-Null
$d.Content.Find.Execute("19×90=1710", $true, $false, $false, $false, $false, $true, 1, $false, "67×31=2077", 2) | Out-Null
$d.Content.Find.Execute("65×65=4225", $true, $false, $false, $false, $false, $true, 1, $false, "72×59=4248", 2) | Out-Null
$d.Content.Find.Execute("25×22=550", $true, $false, $false, $false, $false, $true, 1, $false, "36×91=3276", 2) | Out-Null
$d.Content.Find.Execute("81×95=7695", $true, $false, $false, $false, $false, $true, 1, $false, "49×56=2744", 2) | Out-Null
$d.Content.Find.Execute("55×25=1375", $true, $false, $false, $false, $false, $true, 1, $false, "88×37=3256", 2) | Out-Null
$d.Content.Find.Execute("36×45=1620", $true, $false, $false, $false, $false, $true, 1, $false, "28×51=1428", 2) | Out-Null
$d.Content.Find.Execute("74×84=6216", $true, $false, $false, $false, $false, $true, 1, $false, "51×51=2601", 2) | Out-Null
$d.Content.Find.Execute("30×69=2070", $true, $false, $false, $false, $false, $true, 1, $false, "61×50=3050", 2) | Out-Null
$d.Content.Find.Execute("24×72=1728", $true, $false, $false, $false, $false, $true, 1, $false, "94×52=4888", 2) | Out-Null
$d.Content.Find.Execute("14×16=224", $true, $false, $false, $false, $false, $true, 1, $false, "29×10=290", 2) | Out-Null
$d.Content.Find.Execute("19×98=1862", $true, $false, $false, $false, $false, $true, 1, $false, "32×24=768", 2) | Out-Null
$d.Content.Find.Execute("47×50=2350", $true, $false, $false, $false, $false, $true, 1, $false, "34×41=1394", 2) | Out-Null
$d.Content.Find.Execute("63×58=3654", $true, $false, $false, $false, $false, $true, 1, $false, "25×45=1125", 2) | Out-Null
$d.Content.Find.Execute("41×93=3813", $true, $false, $false, $false, $false, $true, 1, $false, "16×53=848", 2) | Out-Null
$d.Content.Find.Execute("88×78=6864", $true, $false, $false, $false, $false, $true, 1, $false, "58×15=870", 2) | Out-Null
$d.Content.Find.Execute("93×83=7719", $true, $false, $false, $false, $false, $true, 1, $false, "67×13=871", 2) | Out-Null
$d.Content.Find.Execute("36×49=1764", $true, $false, $false, $false, $false, $true, 1, $false, "59×40=2360", 2) | Out-Null
$d.Content.Find.Execute("18×100=1800", $true, $false, $false, $false, $false, $true, 1, $false, "79×35=2765", 2) | Out-Null
$d.Content.Find.Execute("62×91=5642", $true, $false, $false, $false, $false, $true, 1, $false, "56×82=4592", 2) | Out-Null
$d.Content.Find.Execute("93×82=7626", $true, $false, $false, $false, $false, $true, 1, $false, "27×65=1755", 2) | Out-Null
$d.Content.Find.Execute("98×22=2156", $true, $false, $false, $false, $false, $true, 1, $false, "74×13=962", 2) | Out-Null
$d.Content.Find.Execute("71×28=1988", $true, $false, $false, $false, $false, $true, 1, $false, "10×87=870", 2) | Out-Null
$d.Content.Find.Execute("11×76=836", $true, $false, $false, $false, $false, $true, 1, $false, "69×46=3174", 2) | Out-Null
$d.Content.Find.Execute("14×11=154", $true, $false, $false, $false, $false, $true, 1, $false, "97×81=7857", 2) | Out-Null
$d.Content.Find.Execute("53×98=5194", $true, $false, $false, $false, $false, $true, 1, $false, "84×34=2856", 2) | Out-Null
$d.Content.Find.Execute("20×92=1840", $true, $false, $false, $false, $false, $true, 1, $false, "16×64=1024", 2) | Out-Null
$d.Content.Find.Execute("78×80=6240", $true, $false, $false, $false, $false, $true, 1, $false, "79×40=3160", 2) | Out-Null
$d.Content.Find.Execute("53×91=4823", $true, $false, $false, $false, $false, $true, 1, $false, "55×52=2860", 2) | Out-Null
$d.Content.Find.Execute("68×20=1360", $true, $false, $false, $false, $false, $true, 1, $false, "44×100=4400", 2) | Out-Null
$d.Content.Find.Execute("100×88=8800", $true, $false, $false, $false, $false, $true, 1, $false, "76×36=2736", 2) | Out-Null
$d.Content.Find.Execute("39×14=546", $true, $false, $false, $false, $false, $true, 1, $false, "40×61=2440", 2) | Out-Null
$d.Content.Find.Execute("90×70=6300", $true, $false, $false, $false, $false, $true, 1, $false, "35×57=1995", 2) | Out-Null
$d.Content.Find.Execute("26×76=1976", $true, $false, $false, $false, $false, $true, 1, $false, "55×93=5115", 2) | Out-Null
$d.Content.Find.Execute("98×41=4018", $true, $false, $false, $false, $false, $true, 1, $false, "30×57=1710", 2) | Out-Null
$d.Content.Find.Execute("70×78=5460", $true, $false, $false, $false, $false, $true, 1, $false, "71×70=4970", 2) | Out-Null
$d.Content.Find.Execute("98×61=5978", $true, $false, $false, $false, $false, $true, 1, $false, "32×20=640", 2) | Out-Null
$d.Content.Find.Execute("100×76=7600", $true, $false, $false, $false, $false, $true, 1, $false, "62×15=930", 2) | Out-Null
$d.Content.Find.Execute("33×59=1947", $true, $false, $false, $false, $false, $true, 1, $false, "85×45=3825", 2) | Out-Null
$d.Content.Find.Execute("53×66=3498", $true, $false, $false, $false, $false, $true, 1, $false, "76×32=2432", 2) | Out-Null
$d.Content.Find.Execute("89×15=1335", $true, $false, $false, $false, $false, $true, 1, $false, "35×28=980", 2) | Out-Null
$d.Content.Find.Execute("84×48=4032", $true, $false, $false, $false, $false, $true, 1, $false, "100×67=6700", 2) | Out-Null
$d.Content.Find.Execute("14×55=770", $true, $false, $false, $false, $false, $true, 1, $false, "82×94=7708", 2) | Out-Null
$d.Content.Find.Execute("28×79=2212", $true, $false, $false, $false, $false, $true, 1, $false, "42×49=2058", 2) | Out-Null
$d.Content.Find.Execute("61×17=1037", $true, $false, $false, $false, $false, $true, 1, $false, "87×98=8526", 2) | Out-Null
$d.Content.Find.Execute("55×59=3245", $true, $false, $false, $false, $false, $true, 1, $false, "88×31=2728", 2) | Out-Null
$d.Content.Find.Execute("49×51=2499", $true, $false, $false, $false, $false, $true, 1, $false, "30×19=570", 2) | Out-Null
$d.Content.Find.Execute("67×38=2546", $true, $false, $false, $false, $false, $true, 1, $false, "11×66=726", 2) | Out-Null
$d.Content.Find.Execute("65×67=4355", $true, $false, $false, $false, $false, $true, 1, $false, "36×28=1008", 2) | Out-Null
$d.Content.Find.Execute("58×98=5684", $true, $false, $false, $false, $false, $true, 1, $false, "74×70=5180", 2) | Out-Null
$d.Content.Find.Execute("28×20=560", $true, $false, $false, $false, $false, $true, 1, $false, "96×14=1344", 2) | Out-Null
$d.Content.Find.Execute("76×72=5472", $true, $false, $false, $false, $false, $true, 1, $false, "44×30=1320", 2) | Out-Null
$d.Content.Find.Execute("40×90=3600", $true, $false, $false, $false, $false, $true, 1, $false, "22×94=2068", 2) | Out-Null
$d.Content.Find.Execute("21×57=1197", $true, $false, $false, $false, $false, $true, 1, $false, "67×78=5226", 2) | Out-Null
$d.Content.Find.Execute("100×50=5000", $true, $false, $false, $false, $false, $true, 1, $false, "87×41=3567", 2) | Out-Null
$d.Content.Find.Execute("58×24=1392", $true, $false, $false, $false, $false, $true, 1, $false, "71×50=3550", 2) | Out-Null
$d.Content.Find.Execute("94×50=4700", $true, $false, $false, $false, $false, $true, 1, $false, "22×93=2046", 2) | Out-Null
$d.Content.Find.Execute("33×88=2904", $true, $false, $false, $false, $false, $true, 1, $false, "60×81=4860", 2) | Out-Null
$d.Content.Find.Execute("61×65=3965", $true, $false, $false, $false, $false, $true, 1, $false, "53×55=2915", 2) | Out-Null
$d.Content.Find.Execute("92×28=2576", $true, $false, $false, $false, $false, $true, 1, $false, "61×36=2196", 2) | Out-Null
$d.Content.Find.Execute("80×79=6320", $true, $false, $false, $false, $false, $true, 1, $false, "22×86=1892", 2) | Out-Null
$d.Content.Find.Execute("33×77=2541", $true, $false, $false, $false, $false, $true, 1, $false, "61×35=2135", 2) | Out-Null
$d.Content.Find.Execute("69×88=6072", $true, $false, $false, $false, $false, $true, 1, $false, "67×14=938", 2) | Out-Null
$d.Content.Find.Execute("13×80=1040", $true, $false, $false, $false, $false, $true, 1, $false, "88×54=4752", 2) | Out-Null
$d.Content.Find.Execute("43×80=3440", $true, $false, $false, $false, $false, $true, 1, $false, "76×51=3876", 2) | Out-Null
$d.Content.Find.Execute("83×28=2324", $true, $false, $false, $false, $false, $true, 1, $false, "62×85=5270", 2) | Out-Null
$d.Content.Find.Execute("20×67=1340", $true, $false, $false, $false, $false, $true, 1, $false, "10×15=150", 2) | Out-Null
$d.Content.Find.Execute("67×69=4623", $true, $false, $false, $false, $false, $true, 1, $false, "21×22=462", 2) | Out-Null
$d.Content.Find.Execute("49×66=3234", $true, $false, $false, $false, $false, $true, 1, $false, "57×55=3135", 2) | Out-Null
$d.Content.Find.Execute("66×65=4290", $true, $false, $false, $false, $false, $true, 1, $false, "64×80=5120", 2) | Out-Null
$d.Content.Find.Execute("77×29=2233", $true, $false, $false, $false, $false, $true, 1, $false, "72×31=2232", 2) | Out-Null
$d.Content.Find.Execute("18×23=414", $true, $false, $false, $false, $false, $true, 1, $false, "60×92=5520", 2) | Out-Null
$d.Content.Find.Execute("57×19=1083", $true, $false, $false, $false, $false, $true, 1, $false, "52×59=3068", 2) | Out-Null
$d.Content.Find.Execute("99×39=3861", $true, $false, $false, $false, $false, $true, 1, $false, "91×66=6006", 2) | Out-Null
$d.Content.Find.Execute("51×49=2499", $true, $false, $false, $false, $false, $true, 1, $false, "35×34=1190", 2) | Out-Null
$d.Content.Find.Execute("61×99=6039", $true, $false, $false, $false, $false, $true, 1, $false, "36×63=2268", 2) | Out-Null
$d.Content.Find.Execute("17×93=1581", $true, $false, $false, $false, $false, $true, 1, $false, "56×48=2688", 2) | Out-Null
$d.Content.Find.Execute("27×41=1107", $true, $false, $false, $false, $false, $true, 1, $false, "20×11=220", 2) | Out-Null
$d.Content.Find.Execute("83×11=913", $true, $false, $false, $false, $false, $true, 1, $false, "38×94=3572", 2) | Out-Null
$d.Content.Find.Execute("88×41=3608", $true, $false, $false, $false, $false, $true, 1, $false, "72×13=936", 2) | Out-Null
$d.Content.Find.Execute("38×99=3762", $true, $false, $false, $false, $false, $true, 1, $false, "23×93=2139", 2) | Out-Null
$d.Content.Find.Execute("85×64=5440", $true, $false, $false, $false, $false, $true, 1, $false, "85×83=7055", 2) | Out-Null
$d.Content.Find.Execute("67×58=3886", $true, $false, $false, $false, $false, $true, 1, $false, "63×48=3024", 2) | Out-Null
$d.Content.Find.Execute("22×47=1034", $true, $false, $false, $false, $false, $true, 1, $false, "83×76=6308", 2) | Out-Null
$d.Content.Find.Execute("24×35=840", $true, $false, $false, $false, $false, $true, 1, $false, "62×72=4464", 2) | Out-Null
$d.Content.Find.Execute("18×19=342", $true, $false, $false, $false, $false, $true, 1, $false, "51×95=4845", 2) | Out-Null
$d.Content.Find.Execute("47×69=3243", $true, $false, $false, $false, $false, $true, 1, $false, "79×78=6162", 2) | Out-Null
$d.Content.Find.Execute("78×85=6630", $true, $false, $false, $false, $false, $true, 1, $false, "80×11=880", 2) | Out-Null
$d.Content.Find.Execute("68×74=5032", $true, $false, $false, $false, $false, $true, 1, $false, "76×18=1368", 2) | Out-Null
$d.Content.Find.Execute("22×34=748", $true, $false, $false, $false, $false, $true, 1, $false, "68×17=1156", 2) | Out-Null
$d.Content.Find.Execute("98×67=6566", $true, $false, $false, $false, $false, $true, 1, $false, "89×81=7209", 2) | Out-Null
$d.Content.Find.Execute("72×33=2376", $true, $false, $false, $false, $false, $true, 1, $false, "57×88=5016", 2) | Out-Null
$d.Content.Find.Execute("11×27=297", $true, $false, $false, $false, $false, $true, 1, $false, "44×85=3740", 2) | Out-Null
$d.Content.Find.Execute("59×61=3599", $true, $false, $false, $false, $false, $true, 1, $false, "29×92=2668", 2) | Out-Null
$d.Content.Find.Execute("83×89=7387", $true, $false, $false, $false, $false, $true, 1, $false, "27×15=405", 2) | Out-Null
